$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest crypto data.
# NumberFormat '@' forces text entry so numeric-looking prices (e.g. '0.193')
# are not auto-converted to numbers; Style is then reset to 'Normal' so no
# stray formatting/style index is left behind on the cell.
function Set-TextValue($ref, $value) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '69.222.86'
Set-TextValue 'E2' '  +2.04%  '
Set-TextValue 'D3' '3.379.16'
Set-TextValue 'E3' '  +1.35%  '
Set-TextValue 'E4' '  +0.02%  '
Set-TextValue 'D5' '586.35'
Set-TextValue 'E5' '  +0.80%  '
Set-TextValue 'D6' '178.80'
Set-TextValue 'E6' '  +1.24%  '
Set-TextValue 'E7' '  -0.08%  '
Set-TextValue 'E8' '  +1.11%  '
Set-TextValue 'D9' '0.193'
Set-TextValue 'E9' '  +5.54%  '
Set-TextValue 'D10' '0.589'
Set-TextValue 'E10' '  +1.17%  '
Set-TextValue 'D11' '48.32'
Set-TextValue 'E11' '  +2.41%  '
Set-TextValue 'D12' '0.0000280'
Set-TextValue 'E12' '  +2.78%  '
Set-TextValue 'D13' '684.13'
Set-TextValue 'E13' '  -2.02%  '
Set-TextValue 'D14' '8.59'
Set-TextValue 'E14' '  +2.02%  '
Set-TextValue 'D15' '3.924.94'
Set-TextValue 'E15' '  +1.26%  '
Set-TextValue 'D16' '69.251.04'
Set-TextValue 'E16' '  +2.09%  '
Set-TextValue 'D17' '3.379.88'
Set-TextValue 'E17' '  +1.33%  '
Set-TextValue 'D19' '17.61'
Set-TextValue 'E19' '  +0.70%  '
Set-TextValue 'D20' '11.27'
Set-TextValue 'E20' '  +1.91%  '
Set-TextValue 'D21' '0.902'
Set-TextValue 'E21' '  +0.82%  '
Set-TextValue 'D22' '5.42'
Set-TextValue 'E22' '  -0.06%  '
Set-TextValue 'D23' '17.17'
Set-TextValue 'E23' '  +0.89%  '
Set-TextValue 'D24' '103.25'
Set-TextValue 'E24' '  +3.30%  '
Set-TextValue 'D25' '3.92'
Set-TextValue 'E25' '  +0.24%  '
Set-TextValue 'E26' '  +0.69%  '
Set-TextValue 'D27' '9.62'
Set-TextValue 'E27' '  +0.62%  '
Set-TextValue 'D28' '33.87'
Set-TextValue 'E28' '  +2.66%  '
Set-TextValue 'D29' '8.69'
Set-TextValue 'E29' '  +1.43%  '
Set-TextValue 'D30' '6.94'
Set-TextValue 'E30' '  -1.13%  '
Set-TextValue 'D31' '560.57'
Set-TextValue 'E31' '  -1.14%  '
Set-TextValue 'D32' '11.10'
Set-TextValue 'E32' '  +0.93%  '
Set-TextValue 'E33' '  +0.43%  '
Set-TextValue 'D34' '3.55'
Set-TextValue 'E34' '  +5.67%  '
Set-TextValue 'D35' '58.58'
Set-TextValue 'E35' '  +1.89%  '
Set-TextValue 'D36' '0.999'
Set-TextValue 'E36' '  +0.06%  '
Set-TextValue 'D37' '3.679.46'
Set-TextValue 'E37' '  -0.22%  '
Set-TextValue 'D38' '35.73'
Set-TextValue 'E38' '  +2.51%  '
Set-TextValue 'E39' '  +3.94%  '
Set-TextValue 'E40' '  +2.64%  '
Set-TextValue 'D41' '2.67'
Set-TextValue 'E41' '  +1.29%  '
Set-TextValue 'D42' '0.0₃0695'
Set-TextValue 'E42' '  +2.96%  '
Set-TextValue 'E43' '  +0.77%  '
Set-TextValue 'E44' '  +3.49%  '
Set-TextValue 'E45' '  -0.72%  '
Set-TextValue 'E46' '  -0.02%  '
Set-TextValue 'E47' '  +0.94%  '
Set-TextValue 'E48' '  +4.80%  '
Set-TextValue 'E49' '  +0.01%  '
Set-TextValue 'D50' '133.10'
Set-TextValue 'E50' '  +1.60%  '
Set-TextValue 'D51' '2.62'
Set-TextValue 'E51' '  +3.33%  '
